$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "1.00", "555.98") are not silently converted to numbers,
# matching the original inline-string ("General" text) representation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.871.23'
$ws.Range("E2").Value = '  -5.40%  '
$ws.Range("D3").Value = '3.098.51'
$ws.Range("E3").Value = '  -6.08%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '555.98'
$ws.Range("E5").Value = '  -5.25%  '
$ws.Range("D6").Value = '161.38'
$ws.Range("E6").Value = '  -10.46%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.583'
$ws.Range("E8").Value = '  -9.18%  '
$ws.Range("D9").Value = '3.099.10'
$ws.Range("E9").Value = '  -5.96%  '
$ws.Range("D10").Value = '6.71'
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D11").Value = '0.115'
$ws.Range("E11").Value = '  -8.78%  '
$ws.Range("D12").Value = '0.374'
$ws.Range("E12").Value = '  -7.20%  '
$ws.Range("D13").Value = '3.638.24'
$ws.Range("E13").Value = '  -6.15%  '
$ws.Range("E14").Value = '  -2.16%  '
$ws.Range("D15").Value = '63.044.36'
$ws.Range("E15").Value = '  -5.09%  '
$ws.Range("D16").Value = '24.35'
$ws.Range("E16").Value = '  -8.52%  '
$ws.Range("D17").Value = '3.109.04'
$ws.Range("E17").Value = '  -4.80%  '
$ws.Range("E18").Value = '  -7.37%  '
$ws.Range("D19").Value = '399.42'
$ws.Range("E19").Value = '  -6.55%  '
$ws.Range("D20").Value = '12.31'
$ws.Range("E20").Value = '  -5.45%  '
$ws.Range("E21").Value = '  -6.06%  '
$ws.Range("D22").Value = '6.97'
$ws.Range("E22").Value = '  -4.83%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").Value = '67.65'
$ws.Range("E25").Value = '  -5.50%  '
$ws.Range("D26").Value = '0.196'
$ws.Range("E26").Value = '  -4.95%  '
$ws.Range("D27").Value = '0.473'
$ws.Range("E27").Value = '  -8.08%  '
$ws.Range("D28").Value = '0.0₃0994'
$ws.Range("E28").Value = '  -13.12%  '
$ws.Range("D29").Value = '8.57'
$ws.Range("E29").Value = '  -5.95%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("E32").Value = '  -8.07%  '
$ws.Range("D33").Value = '20.77'
$ws.Range("E33").Value = '  -7.00%  '
$ws.Range("D34").Value = '4.76'
$ws.Range("E34").Value = '  -7.75%  '
$ws.Range("D35").Value = '6.12'
$ws.Range("E35").Value = '  -6.75%  '
$ws.Range("D36").Value = '151.49'
$ws.Range("E36").Value = '  -4.64%  '
$ws.Range("E37").Value = '  -8.56%  '
$ws.Range("E38").Value = '  -8.57%  '
$ws.Range("D39").Value = '2.687.85'
$ws.Range("E39").Value = '  -6.46%  '
$ws.Range("E40").Value = '  -9.54%  '
$ws.Range("D41").Value = '23.18'
$ws.Range("E41").Value = '  -11.71%  '
$ws.Range("D42").Value = '3.98'
$ws.Range("E42").Value = '  -7.73%  '
$ws.Range("D43").Value = '38.11'
$ws.Range("E43").Value = '  -3.97%  '
$ws.Range("D44").Value = '0.691'
$ws.Range("E44").Value = '  -8.00%  '
$ws.Range("D45").Value = '0.0598'
$ws.Range("E45").Value = '  -8.89%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0251'
$ws.Range("E46").Value = '  -7.47%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '5.14'
$ws.Range("E47").Value = '  -13.06%  '
$ws.Range("D48").Value = '281.23'
$ws.Range("E48").Value = '  -9.55%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = '20.47'
$ws.Range("E50").Value = '  -10.68%  '
$ws.Range("D51").Value = '0.0962'
$ws.Range("E51").Value = '  -6.33%  '
